# edit.ps1
# Applies the "small improvements in questions.xls" commit:
#  1. Workbook view: clear the saved activeTab marker (bookViews/workbookView).
#  2. Sheet1 ("Лист1"): add a "Фамилия Имя Отчество" (respondent full name) value
#     in column A for rows 2-120, taken from the new shared-string list.
#  3. Sheet1: clear leftover scratch numbers in B163:D165 and C166:C168.
#  4. Sheet1 view: move the active selection / top-left cell, mark sheet1 as
#     the tab that is now selected.
#  5. Sheet2 ("Лист2") view: it is no longer the selected tab; selection moves
#     to AC17.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Лист1")
$sheet2 = $wb.Worksheets.Item("Лист2")

$names = @(
    'Александрова-Зорина Елизавета Борисовна',
    'Аншаков Михаил Геннадьевич',
    'Ашурков Владимир Львович',
    'Бадыкова Фиалка Ахмадеевна',
    'Бакиров Игорь Вакильевич',
    'Баласанов Андрей Евгеньевич',
    'Барабаш Кирилл Владимирович',
    'Баронова Мария Николаевна',
    'Беззуб Алексей Юрьевич',
    'Безруких Олег Анатольевич',
    'Билунов Денис Борисович',
    'Блиндул Алексей Валерьевич',
    'Брусиловский Максим Анатольевич',
    'Будникова Софья Владимировна',
    'Быстров Андрей Сергеевич',
    'Васильева Елена Борисовна',
    'Виноградов Сергей Николаевич',
    'Винокуров Александр Иванович',
    'Витухновская Алина Александровна',
    'Власов Сергей Игоревич',
    'Гаврилов Андрей Игоревич',
    'Газарян Сурен Владимирович',
    'Галямина Юлия Евгеньевна',
    'Гарначук Владимир Федорович',
    'Гельфанд Михаил Сергеевич',
    'Глускин Владимир Александрович',
    'Головин Дмитрий Александрович',
    'Гонгальский Максим Брониславович',
    'Горник Александр Львович',
    'Готсданкер Алексей Сергеевич',
    'Гребнева Ирина Георгиевна',
    'Давыденко Денис Вячеславович',
    'Давыдов Андрей Владимирович',
    'Дегтярь Иван Васильевич',
    'Демидов Михаил Александрович',
    'Дергачев Вадим Александрович',
    'Дзядко Филипп Викторович',
    'Доможиров Евгений Валерьевич',
    'Езеев Федор Андреевич',
    'Зорин Константин Игоревич',
    'Иванов Андрей Геннадьевич',
    'Илларионов Андрей Николаевич',
    'Кара-Мурза Владимир Владимирович',
    'Каржаева Неонила Васильевна',
    'Каспаров Гарри Кимович',
    'Кац Максим Евгеньевич',
    'Кашин Олег Владимирович',
    'Козырев Олег Вилисович',
    'Колчинцев Вадим Валерьевич',
    'Крашенинников Федор Геннадиевич',
    'Крылов Олег Анатольевич',
    'Крюков Василий Анатольевич',
    'Кузин Евгений Андреевич',
    'Кузнецов Андрей Владимирович',
    'Курамшин Владимир Вячеславович',
    'Лавров Андрей Валерьевич',
    'Лазуренко (Северский) Артем Сергеевич',
    'Левченко Екатерина Валентиновна',
    'Левшиц Николай Дмитриевич',
    'Литвинов Георгий Александрович (Артём Драгунов)',
    'Магкоева Белла (Изабель) Казбековна',
    'Малышев Владимир Эдуардович',
    'Мальцева Анастасия Анатольевна (Анастасия Хрустальная)',
    'Матвеев Михаил Николаевич',
    'Мирзоев Владимир Владимирович',
    'Митюшкина Надежда Львовна',
    'Мокшанов Александр Александрович',
    'Мухин Юрий Игнатьевич',
    'Некрасов Дмитрий Александрович',
    'Образцова Алиса Сергеевна',
    'Овдиенко Игорь Геннадьевич',
    'Ольшанский Леонид Дмитриевич',
    'Осенин Владимир Олегович',
    'Отставных Валерий Владимирович',
    'Пархоменко Сергей Борисович',
    'Первушин Александр Сергеевич',
    'Петречук Лариса Леонидовна',
    'Пионтковский Андрей Андреевич',
    'Поляков Анатолий Викторович',
    'Пономарев Илья Владимирович',
    'Поткин (Басманов) Владимир Анатольевич',
    'Пряников Павел Николаевич',
    'Ренёв Денис Владимирович',
    'Русакова Елена Леонидовна',
    'Савостин Михаил Олегович',
    'Сайдашев Радик Ромович',
    'Семенов Владимир Матвеевич',
    'Семенов Игорь Вячеславович',
    'Скалаух Иван Сергеевич',
    'Смирнов Сергей Сергеевич',
    'Спорыхина Ульяна Викторовна',
    'Стефанов Борис Александрович',
    'Сухарева Татьяна Викторовна',
    'Терегулов Артур Ринатович',
    'Удальцова Анастасия Олеговна',
    'Чупров Алексей Геннадьевич',
    'Шатов Станислав Николаевич',
    'Шац Михаил Григорьевич',
    'Шнейдер Михаил Яковлевич',
    'Щербаков Александр Вениаминович',
    'Эсауленко Дмитрий Николаевич',
    'Янкаускас Константин Стасисович',
    'Яшин (Ясин) Игорь Геннадьевич',
    'Аитова Екатерина Петровна',
    'Волкова Александра Ивановна (Женя Отто)',
    'Николаев Александр Александрович',
    'Печенев Александр Сергеевич',
    'Санников Максим Андреевич',
    'Шалимов Роман Николаевич',
    'Давидис Сергей Константинович',
    'Долгих Антон Витальевич',
    'Залесский Александр Валерьевич',
    'Пивоваров Андрей Сергеевич',
    'Тютрин Иван Иванович',
    'Шальнев Андрей Сергеевич',
    'Артёмов Игорь Владимирович',
    'Демушкин Дмитрий Николаевич',
    'Дровецкий Василий Валерьевич',
    'Крылов Константин Анатольевич'
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $sheet1.Cells.Item($row, 1).Value = $names[$i]
}

# Clear the stale scratch values left in this little side table.
$sheet1.Range("B163:D165").ClearContents()
$sheet1.Range("C166:C168").ClearContents()

# Restore the view state captured in the saved workbook: Лист1 becomes the
# active/selected sheet, scrolled down with a new selection.
[void]$sheet1.Activate()
$excel.ActiveWindow.ScrollRow = 136
[void]$sheet1.Range("B162:G175").Select()

# Лист2 keeps its own remembered selection, but it is no longer the tab that
# is marked selected when the file is reopened.
[void]$sheet2.Range("AC17").Select()

# Re-activate Лист1 last so it is the workbook's active sheet on save.
[void]$sheet1.Activate()
